# Applies the "Corrected spelling mistake and removed unecessary fields
# used for testing" commit to the Notification List workbook.
#
# Changes made on the "Templates" sheet:
#   1. Fix the "queruies" -> "queries" spelling mistake in the two
#      Event Template Example cells (B2 and C2).
#   2. Remove the unnecessary/duplicate testing row (row 3: "Ubaya" /
#      the "0.0" templates) that was only used for testing.
#   3. Reflow the "First Message" template paragraph (now row 3 after
#      the deletion above) onto multiple lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# --- 1. Fix the spelling mistake in the two template examples (row 2) ---
$eventTypeTemplate = $ws.Range("B2").Value()
$eventTypeTemplate = $eventTypeTemplate.Replace("queruies", "queries")
$ws.Range("B2").Value = $eventTypeTemplate

$ubayamTemplate = $ws.Range("C2").Value()
$ubayamTemplate = $ubayamTemplate.Replace("queruies", "queries")
$ws.Range("C2").Value = $ubayamTemplate

# --- 2. Remove the unnecessary testing row ("Ubaya" / "0.0" templates) ---
$ws.Rows.Item(3).Delete()

# --- 3. Reflow the first-message paragraph (now row 3) onto more lines ---
$oldParagraph = "We have embarked on a new initiative to inform you via WhatsApp on key SSPT specific events and your ubayams & services at the temple. This is one our efforts to improve our engagement with you. We seek you understanding while we finetune the initiative. Looking forward to your support and constructive feedback."
$newParagraph = "We have embarked on a new initiative to inform you via WhatsApp on key SSPT specific events and your ubayams & services at the temple. `nThis is one our efforts to improve our engagement with you. We seek you understanding while we finetune the initiative. `nLooking forward to your support and constructive feedback."

$firstMessage = $ws.Range("B3").Value()
$firstMessage = $firstMessage.Replace($oldParagraph, $newParagraph)
$ws.Range("B3").Value = $firstMessage

# Row 3 now holds noticeably more text/lines, so grow it to fit.
$ws.Rows.Item(3).RowHeight = 238

# Reflect the new "current cell" the author ended up on after the edits.
$ws.Range("C3").Select()
